# Auto-generated edit script applying numeric corrections to the
# Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values were recomputed upstream (scheduled market-data refresh); this
# script writes the refreshed currentAveragePrice* / LevePrice* /
# LeveProfit* figures (columns H:N) for the affected leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 20835208
$ws.Range("I112").Value = 2605
$ws.Range("J112").Value = 41667810
$ws.Range("K112").Value = 7815
$ws.Range("L112").Value = 125003430
$ws.Range("M112").Value = -6707
$ws.Range("N112").Value = -125005646
$ws.Range("H115").Value = 1238.5714
$ws.Range("I115").Value = 504
$ws.Range("J115").Value = 3075
$ws.Range("K115").Value = 1512
$ws.Range("L115").Value = 9225
$ws.Range("M115").Value = 55
$ws.Range("N115").Value = -12359
$ws.Range("H118").Value = 1750.8064
$ws.Range("I118").Value = 615
$ws.Range("J118").Value = 2023.4
$ws.Range("K118").Value = 1845
$ws.Range("L118").Value = 6070.200000000001
$ws.Range("M118").Value = -188
$ws.Range("N118").Value = -9384.200000000001
$ws.Range("H127").Value = 1241979
$ws.Range("I127").Value = 3565
$ws.Range("J127").Value = 1403511.2
$ws.Range("K127").Value = 10695
$ws.Range("L127").Value = 4210533.6
$ws.Range("M127").Value = -5735
$ws.Range("N127").Value = -4220453.6
$ws.Range("H129").Value = 4546519.5
$ws.Range("I129").Value = 125003000
$ws.Range("J129").Value = 992.1698
$ws.Range("K129").Value = 375009000
$ws.Range("L129").Value = 2976.5094
$ws.Range("M129").Value = -375004000
$ws.Range("N129").Value = -12976.5094
$ws.Range("H137").Value = 2329363.8
$ws.Range("I137").Value = 4550598.5
$ws.Range("J137").Value = 2355.762
$ws.Range("K137").Value = 13651795.5
$ws.Range("L137").Value = 7067.286
$ws.Range("M137").Value = -13649245.5
$ws.Range("N137").Value = -12167.286
$ws.Range("H138").Value = 4617.229
$ws.Range("I138").Value = 3050.2307
$ws.Range("J138").Value = 5199.2573
$ws.Range("K138").Value = 9150.6921
$ws.Range("L138").Value = 15597.7719
$ws.Range("M138").Value = -4010.6921
$ws.Range("N138").Value = -25877.7719

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3480
$ws.Range("I61").Value = 2520
$ws.Range("K61").Value = 2520
$ws.Range("M61").Value = -2308
$ws.Range("H74").Value = 3235
$ws.Range("I74").Value = 1904
$ws.Range("J74").Value = 3900.5
$ws.Range("K74").Value = 1904
$ws.Range("L74").Value = 3900.5
$ws.Range("M74").Value = -1030
$ws.Range("N74").Value = -5648.5
$ws.Range("H77").Value = 3235
$ws.Range("I77").Value = 1904
$ws.Range("J77").Value = 3900.5
$ws.Range("K77").Value = 9520
$ws.Range("L77").Value = 19502.5
$ws.Range("M77").Value = -5152
$ws.Range("N77").Value = -28238.5
$ws.Range("H112").Value = 27462.334
$ws.Range("J112").Value = 27462.334
$ws.Range("L112").Value = 27462.334
$ws.Range("N112").Value = -30416.334
$ws.Range("H122").Value = 4563.636
$ws.Range("I122").Value = 3033.3333
$ws.Range("J122").Value = 6400
$ws.Range("K122").Value = 9099.999899999999
$ws.Range("L122").Value = 19200
$ws.Range("M122").Value = -6649.999899999999
$ws.Range("N122").Value = -24100
$ws.Range("H132").Value = 3575.1155
$ws.Range("I132").Value = 3155.8
$ws.Range("J132").Value = 3837.1875
$ws.Range("K132").Value = 9467.400000000001
$ws.Range("L132").Value = 11511.5625
$ws.Range("M132").Value = -6937.400000000001
$ws.Range("N132").Value = -16571.5625
$ws.Range("H136").Value = 3480
$ws.Range("I136").Value = 2520
$ws.Range("K136").Value = 7560
$ws.Range("M136").Value = -5010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 883.85
$ws.Range("I94").Value = 743.86664
$ws.Range("J94").Value = 1303.8
$ws.Range("K94").Value = 743.86664
$ws.Range("L94").Value = 1303.8
$ws.Range("M94").Value = -292.86664
$ws.Range("N94").Value = -2205.8
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H134").Value = 3195.625
$ws.Range("I134").Value = 3091.2646
$ws.Range("K134").Value = 9273.793799999999
$ws.Range("M134").Value = -6738.793799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 658053.9
$ws.Range("I6").Value = 1502631.8
$ws.Range("J6").Value = 95002
$ws.Range("K6").Value = 1502631.8
$ws.Range("L6").Value = 95002
$ws.Range("M6").Value = -1502518.8
$ws.Range("N6").Value = -95228
$ws.Range("H7").Value = 80
$ws.Range("I7").Value = 80
$ws.Range("K7").Value = 80
$ws.Range("M7").Value = 33
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null
$ws.Range("H31").Value = 4353591
$ws.Range("I31").Value = 10004382
$ws.Range("J31").Value = 6828.769
$ws.Range("K31").Value = 10004382
$ws.Range("L31").Value = 6828.769
$ws.Range("M31").Value = -10004087
$ws.Range("N31").Value = -7418.769
$ws.Range("H34").Value = 4353591
$ws.Range("I34").Value = 10004382
$ws.Range("J34").Value = 6828.769
$ws.Range("K34").Value = 10004382
$ws.Range("L34").Value = 6828.769
$ws.Range("M34").Value = -10004180
$ws.Range("N34").Value = -7232.769
$ws.Range("H58").Value = 10873240
$ws.Range("I58").Value = 1853.7391
$ws.Range("K58").Value = 1853.7391
$ws.Range("M58").Value = -1650.7391
$ws.Range("H123").Value = 30926.666
$ws.Range("J123").Value = 30926.666
$ws.Range("L123").Value = 30926.666
$ws.Range("N123").Value = -40726.666
$ws.Range("H132").Value = 3764.875
$ws.Range("I132").Value = 3653.6
$ws.Range("J132").Value = 3950.3333
$ws.Range("K132").Value = 10960.8
$ws.Range("L132").Value = 11850.9999
$ws.Range("M132").Value = -8430.799999999999
$ws.Range("N132").Value = -16910.9999
$ws.Range("H134").Value = 6068.4
$ws.Range("I134").Value = 4556.8887
$ws.Range("J134").Value = 8335.666999999999
$ws.Range("K134").Value = 13670.6661
$ws.Range("L134").Value = 25007.001
$ws.Range("M134").Value = -11135.6661
$ws.Range("N134").Value = -30077.001
$ws.Range("H136").Value = 10873240
$ws.Range("I136").Value = 1853.7391
$ws.Range("K136").Value = 5561.2173
$ws.Range("M136").Value = -3011.2173
$ws.Range("H140").Value = 59225
$ws.Range("J140").Value = 59225
$ws.Range("L140").Value = 59225
$ws.Range("N140").Value = -69585

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2564925
$ws.Range("J113").Value = 887.46875
$ws.Range("L113").Value = 2662.40625
$ws.Range("N113").Value = -7002.40625
$ws.Range("H131").Value = 1486.8704
$ws.Range("J131").Value = 1120.9807
$ws.Range("L131").Value = 3362.9421
$ws.Range("N131").Value = -13442.9421

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 47970.8
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 59513.5
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 59513.5
$ws.Range("M40").Value = -1649
$ws.Range("N40").Value = -59815.5
$ws.Range("H46").Value = 13940.571
$ws.Range("I46").Value = 9000
$ws.Range("K46").Value = 9000
$ws.Range("M46").Value = -8844
$ws.Range("H80").Value = 2915.0833
$ws.Range("I80").Value = 3112.4614
$ws.Range("J80").Value = 2681.818
$ws.Range("K80").Value = 3112.4614
$ws.Range("L80").Value = 2681.818
$ws.Range("M80").Value = -2114.4614
$ws.Range("N80").Value = -4677.818
$ws.Range("H83").Value = 2915.0833
$ws.Range("I83").Value = 3112.4614
$ws.Range("J83").Value = 2681.818
$ws.Range("K83").Value = 15562.307
$ws.Range("L83").Value = 13409.09
$ws.Range("M83").Value = -10570.307
$ws.Range("N83").Value = -23393.09
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352
$ws.Range("H113").Value = 3800
$ws.Range("I113").Value = 3120
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 3120
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -950
$ws.Range("N113").Value = -9840
$ws.Range("H122").Value = 7933.778
$ws.Range("I122").Value = 6133.3335
$ws.Range("J122").Value = 8834
$ws.Range("K122").Value = 18400.0005
$ws.Range("L122").Value = 26502
$ws.Range("M122").Value = -15950.0005
$ws.Range("N122").Value = -31402
$ws.Range("H126").Value = 592393.6
$ws.Range("I126").Value = 2555.8
$ws.Range("J126").Value = 838159.4399999999
$ws.Range("K126").Value = 7667.400000000001
$ws.Range("L126").Value = 2514478.32
$ws.Range("M126").Value = -5197.400000000001
$ws.Range("N126").Value = -2519418.32
$ws.Range("H132").Value = 4638.0557
$ws.Range("I132").Value = 5364.3687
$ws.Range("J132").Value = 3826.2942
$ws.Range("K132").Value = 16093.1061
$ws.Range("L132").Value = 11478.8826
$ws.Range("M132").Value = -13563.1061
$ws.Range("N132").Value = -16538.8826

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 16666.666
$ws.Range("J63").Value = 16666.666
$ws.Range("L63").Value = 16666.666
$ws.Range("N63").Value = -18164.666
$ws.Range("H66").Value = 16666.666
$ws.Range("J66").Value = 16666.666
$ws.Range("L66").Value = 49999.99800000001
$ws.Range("N66").Value = -57487.99800000001
$ws.Range("H110").Value = 30322
$ws.Range("J110").Value = 30322
$ws.Range("L110").Value = 30322
$ws.Range("N110").Value = -38502
$ws.Range("H136").Value = 2443644.2
$ws.Range("I136").Value = 3849554
$ws.Range("K136").Value = 11548662
$ws.Range("M136").Value = -11546112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 30000.666
$ws.Range("J129").Value = 30000.666
$ws.Range("L129").Value = 30000.666
$ws.Range("N129").Value = -40000.666
$ws.Range("H132").Value = 1547354.6
$ws.Range("I132").Value = 1963589
$ws.Range("J132").Value = 31072
$ws.Range("K132").Value = 5890767
$ws.Range("L132").Value = 93216
$ws.Range("M132").Value = -5888237
$ws.Range("N132").Value = -98276
$ws.Range("H133").Value = 36400
$ws.Range("J133").Value = 36400
$ws.Range("L133").Value = 36400
$ws.Range("N133").Value = -46520
$ws.Range("H138").Value = 32164.5
$ws.Range("J138").Value = 32164.5
$ws.Range("L138").Value = 32164.5
$ws.Range("N138").Value = -42444.5

